$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
$data = @{
    2 = @(1.02, 1.014578636771362, 1.02109545282962, 0.9926147277508489, 1.026734190128177, 1.027356169029586, 1.019809068733895, 1.023933795138616, 0.9955398523336033, 1.02955597232027, 1.010732584082031)
    3 = @(1.02, 1.015372825508716, 1.021648998064177, 0.9936372048519304, 1.027775913231664, 1.027454466200707, 1.020238716758754, 1.024294413092469, 0.9963617723202692, 1.030404638320952, 1.010875357983942)
    4 = @(1.02, 1.01588719143872, 1.022007399451759, 0.9942998659930995, 1.028450652445256, 1.027516807935385, 1.020516557040173, 1.024527272760225, 0.9968940712668345, 1.030953859791596, 1.010967657312384)
    5 = @(1.02, 1.016103542540507, 1.02215812252836, 0.9945786998346017, 1.028734473301778, 1.027542713420128, 1.020633319147839, 1.024625049900427, 0.997117960005301, 1.031184770103097, 1.011006439227543)
    6 = @(1.02, 1.016139875346889, 1.022183432545663, 0.9946255319796338, 1.028782137453733, 1.027547045280939, 1.020652921528421, 1.024641460214893, 0.9971555583673453, 1.031223541963103, 1.011012949655234)
    7 = @(1.02, 1.015890081895846, 1.022009413222087, 0.9943035907982488, 1.028454444244105, 1.027517155276983, 1.020518117387344, 1.024528579724563, 0.9968970624462087, 1.030956945158587, 1.010968175600148)
    8 = @(1.02, 1.014846937791584, 1.021282479346451, 0.9929600610674301, 1.027086105002159, 1.027389650198935, 1.019954304997863, 1.024055767254992, 0.995817528259106, 1.029842766238608, 1.010780852542736)
    9 = @(1.02, 1.013012471430723, 1.0200032926759, 0.9906006454969559, 1.024680129559854, 1.027155322121942, 1.01895953549275, 1.023218952291585, 0.9939188001724441, 1.027880081850217, 1.010450132069122)
    10 = @(1.02, 1.011792061781477, 1.019151787259163, 0.989033133672735, 1.023079711057067, 1.026992651759648, 1.018295572401836, 1.022658682655748, 0.9926553831429383, 1.026572116476594, 1.01022924902181)
    11 = @(1.02, 1.011264238655276, 1.018783402029491, 0.988355674866747, 1.022387569971828, 1.026920690783925, 1.018007896601439, 1.022415526229124, 0.9921088820399291, 1.026005881666279, 1.010133513541084)
    12 = @(1.02, 1.011068276624951, 1.018646617466701, 0.9881042295826724, 1.022130606686016, 1.026893732933828, 1.017901015712928, 1.022325124700735, 0.9919059725120875, 1.025795576230091, 1.010097939735453)
    13 = @(1.02, 1.011110306823752, 1.018675955924722, 0.9881581567098651, 1.022185720342872, 1.026899525807891, 1.017923943153565, 1.022344519846423, 0.9919494934313052, 1.025840686567911, 1.010105571037052)
    14 = @(1.02, 1.011248038426557, 1.018772094342731, 0.9883348863814464, 1.022366326675028, 1.026918467094392, 1.017999062303558, 1.022408055285129, 0.9920921077337197, 1.025988497352497, 1.010130573270043)
    15 = @(1.02, 1.011332912043005, 1.018831335086416, 0.9884438009545853, 1.022477621282854, 1.026930107205963, 1.018045342372096, 1.022447190679979, 0.9921799884222134, 1.026079571051711, 1.010145976208795)
    16 = @(1.02, 1.011827104902073, 1.019176242708107, 0.9890781214508737, 1.023125664183213, 1.026997395519953, 1.018314660892675, 1.022674808529186, 0.9926916645766087, 1.02660969828329, 1.010235600772084)
    17 = @(1.02, 1.012137266683486, 1.019392681621267, 0.989476357848556, 1.023532392778205, 1.027039196182944, 1.018483551077508, 1.02281743913593, 0.9930127773699352, 1.026942266630235, 1.010291795624537)
    18 = @(1.02, 1.012318238746825, 1.019518957754618, 0.9897087662937556, 1.023769712703124, 1.027063430755536, 1.018582044870008, 1.022900579436725, 0.9932001317071769, 1.027136259908895, 1.010324564261182)
    19 = @(1.02, 1.012379955666002, 1.019562019855237, 0.9897880325774034, 1.023850646548029, 1.027071669159487, 1.018615625792904, 1.022928918984847, 0.9932640239640975, 1.027202408574961, 1.010335736003192)
    20 = @(1.02, 1.012103983048682, 1.019369456557838, 0.9894336180360679, 1.023488746159093, 1.027034726573043, 1.018465432509207, 1.022802141769954, 0.9929783193494215, 1.0269065839679, 1.010285767363373)
    21 = @(1.02, 1.01120747726855, 1.018743782577787, 0.9882828385668249, 1.022313139050654, 1.026912895656335, 1.017976942275697, 1.022389347947773, 0.9920501090198102, 1.025944970217941, 1.010123211102001)
    22 = @(1.02, 1.010644358764085, 1.018350687295426, 0.9875604150241495, 1.021574733459806, 1.026834974902828, 1.017669663298203, 1.022129331873025, 0.9914670000341481, 1.025340478099863, 1.010020928274945)
    23 = @(1.02, 1.010942825882995, 1.01855904641264, 0.9879432794643023, 1.021966105282852, 1.026876407176914, 1.017832571161205, 1.022267216078773, 0.991776070289318, 1.02566091982303, 1.010075157530907)
    24 = @(1.02, 1.012119022305524, 1.01937995086879, 0.9894529299347244, 1.023508467934137, 1.027036746651779, 1.018473619562951, 1.022809054157287, 0.9929938892766442, 1.026922707390686, 1.010288491302562)
    25 = @(1.02, 1.013486278870277, 1.020333773773821, 0.9912096547607049, 1.025301507920542, 1.027217041258083, 1.019216850000661, 1.023435715704151, 0.9944092447426414, 1.028387400674444, 1.010535703909455)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
